$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the existing data (header row 1 + data rows 2..27)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$firstDataRow = 2
$lastDataRow = $lastRow

$numDataRows = $lastDataRow - $firstDataRow + 1

# Source range containing the original data rows (ids 1..26)
$srcRange = $ws.Range("A$($firstDataRow):F$($lastDataRow)")

# Destination starts right after the current last row
$destStartRow = $lastDataRow + 1
$destRange = $ws.Range("A$($destStartRow):F$($destStartRow + $numDataRows - 1)")

# Copy the original rows into the new rows below
$srcRange.Copy($destRange)

# Fix up the id column (A) in the newly appended rows so ids continue
# incrementing (27, 28, 29, ...) instead of repeating (1, 2, 3, ...)
for ($i = 0; $i -lt $numDataRows; $i++) {
    $srcRow = $firstDataRow + $i
    $destRow = $destStartRow + $i
    $idValue = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($destRow, 1).Value = $idValue + $numDataRows
}

$excel.CutCopyMode = $false
